# Generate Report for Archive
#
# 1. Status text update: "Ready for handoff" -> "In Translation"
#    (appears on the Overview sheet in columns E/F row 2, and on the
#    zh-cn / de-de sheets in column C row 2 - these are the per-language
#    "Status" cells for the single tracked file.)
#
# 2. Column width update: the "Status" columns (Overview!E:F, zh-cn!C,
#    de-de!C) are narrowed. ColumnWidth is quantized by Excel to 1/6 of a
#    character, so 12.5 is the closest settable value that reproduces the
#    target width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column widths ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
